$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Row 5: fix the /games filter endpoint text (missing "=")
$ws.Range("D5").Value = "/games?result={WHITE/BLACK/DRAW}"

# Row 6: mark as done instead of GET/PATCH
$ws.Range("J6").Value = "done"

# Row 8 (/players): PUT/PATCH columns no longer apply (moved down to /players/{id}), just GET/POST now
$ws.Range("G8").Value = "x"
$ws.Range("H8").Value = "x"
$ws.Range("I8").Value = "x"
$ws.Range("J8").Value = "done"

# Row 9 (/players/{id}): add PUT (update players data) and PATCH (delete player... actually update/delete)
$ws.Range("G9").Value = "update players data"
$ws.Range("H9").Value = "delete player"
$ws.Range("J9").Value = "GET/PUT/PATCH"
